$p = $ppt.ActivePresentation

# Setting TextRange.Text to a value whose concatenated characters already
# match the current text (just split across multiple runs) is treated as a
# no-op by this host and the run merge never reaches the saved XML. Routing
# the assignment through a throwaway placeholder value first forces the
# engine to actually rebuild the paragraph into a single run.
function Set-MergedText($shape, [string]$text) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "___tmp___"
    $tr.Text = $text
}

# --- Slide 11: "Saisir un élève" / "Saisir le prénom d'un élève" ------------
# The original runs were split ("Saisir " + "un élève"); merge them into a
# single run per the target deck.
$s11 = $p.Slides.Item(11)
Set-MergedText $s11.Shapes.Item("ZoneTexte 5") "Saisir un élève"
Set-MergedText $s11.Shapes.Item("ZoneTexte 6") "Saisir le prénom d’un élève"

# --- Slide 13: "Storyboard - Etudiants" title --------------------------------
# Merge the split " - " and "Etudiants" runs into the trailing run's text,
# leaving the separate leading "Storyboard" run (and its err="1" rPr) alone.
$s13 = $p.Slides.Item(13)
$title13 = $s13.Shapes.Item("ZoneTexte 95").TextFrame.TextRange
$prefixLen = 10 # length of "Storyboard"
$tail = $title13.Characters($prefixLen + 1, $title13.Length - $prefixLen)
$tail.Text = "___tmp___"
$title13b = $s13.Shapes.Item("ZoneTexte 95").TextFrame.TextRange
$tail2 = $title13b.Characters($prefixLen + 1, $title13b.Length - $prefixLen)
$tail2.Text = " - Etudiants"

# --- Slide 3: login button "OK" -> "Connexion" (student login, grouped) -----
$s3 = $p.Slides.Item(3)
$grp1 = $s3.Shapes.Item("Groupe 1")
$grp1.GroupItems.Item("Rectangle 9").TextFrame.TextRange.Text = "Connexion"

# --- Slide 4: login button "OK" -> "Connexion" (admin login) ----------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item("Rectangle 9").TextFrame.TextRange.Text = "Connexion"
